# Append a new data row (row 9) to Sheet1, mirroring the existing table:
#   A9 = 7                          (numeric id, like A2:A8)
#   B9 = "test test test 4 10"      (comment text)
#   C9 = "04-10-2023"               (date kept as literal text, like C2:C8)
#
# Column C's existing rows store the date as text (not a real date value),
# so the new cell must stay text too. Excel would otherwise auto-convert a
# "mm-dd-yyyy"-looking string into a date serial number, so we force the
# cell to Text format first, assign the literal value, then drop the cell
# back to the default "Normal" style (matching the rest of the column,
# which carries no explicit style) to avoid leaving a stray number format
# on just this one cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A9").Value = 7
$ws.Range("B9").Value = "test test test 4 10"

$ws.Range("C9").NumberFormat = "@"
$ws.Range("C9").Value = "04-10-2023"
$ws.Range("C9").Style = "Normal"
